# Generate Report for Handoff
#
# The af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.md file moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with refreshed
# handoff timestamps and a new error detail message recorded on the
# per-language sheets. Update all three sheets accordingly.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef4c4f86babbf720b66eec4f1b1dec20daac7375/e2e/af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d41ed0725b3168b6c29c8aa912460f46e36139af/e2e/af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.md."

# --- Overview sheet: row 3 is the af3f0b90 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = "2016-08-30 19:00:58"

# --- zh-cn sheet: row 3 is the af3f0b90 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = "2016-08-30 19:00:54"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.142857142857146

# --- de-de sheet: row 3 is the af3f0b90 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = "2016-08-30 19:00:58"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.142857142857146
